$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 09:24"

# --- Country name swaps (order of appearance in the data changes) ---
# Emiratos Arabes Unidos / Ucrania swap places (row 38 <-> row 39)
$ws.Range("A38").Value = "Ucrania"
$ws.Range("A39").Value = "Emiratos Arabes Unidos"

# Suiza / Armenia swap places (row 52 <-> row 53)
$ws.Range("A52").Value = "Armenia"
$ws.Range("A53").Value = "Suiza"

# Groenlandia / Islas Malvinas swap places (row 209 <-> row 210)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Numeric data refresh ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3617040
$ws.Range("C4").Value = 213
$ws.Range("D4").Value = 1645966
$ws.Range("E4").Value = 1830924
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 140150

# Row 34 - Belgica
$ws.Range("B34").Value = 63039
$ws.Range("C34").Value = 167
$ws.Range("D34").Value = 17242
$ws.Range("E34").Value = 36005
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 9792

# Row 38 - now Ucrania
$ws.Range("B38").Value = 56455
$ws.Range("C38").Value = 848
$ws.Range("D38").Value = 28931
$ws.Range("E38").Value = 26079
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 1445

# Row 39 - now Emiratos Arabes Unidos
$ws.Range("B39").Value = 55848
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 46418
$ws.Range("E39").Value = 9095
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 335

# Row 48 - Afganistan
$ws.Range("B48").Value = 35070
$ws.Range("C48").Value = 76
$ws.Range("D48").Value = 22824
$ws.Range("E48").Value = 11133
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 19
$ws.Range("H48").Value = 1113

# Row 52 - now Armenia
$ws.Range("B52").Value = 33559
$ws.Range("C52").Value = 554
$ws.Range("D52").Value = 21931
$ws.Range("E52").Value = 11021
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 15
$ws.Range("H52").Value = 607

# Row 53 - now Suiza
$ws.Range("B53").Value = 33148
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 29800
$ws.Range("E53").Value = 1380
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1968

# Row 61 - Moldavia
$ws.Range("D61").Value = 13640
$ws.Range("E61").Value = 5738
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 662

# Row 75 - El Salvador
$ws.Range("D75").Value = 6128
$ws.Range("E75").Value = 4219
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 12
$ws.Range("H75").Value = 298

# Row 99 - Hungria
$ws.Range("B99").Value = 4279
$ws.Range("C99").Value = 16
$ws.Range("D99").Value = 3156
$ws.Range("E99").Value = 528

# Row 110 - Sri Lanka
$ws.Range("B110").Value = 2674
$ws.Range("C110").Value = 3
$ws.Range("D110").Value = 2001
$ws.Range("E110").Value = 662

# Row 137 - Letonia
$ws.Range("B137").Value = 1179
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 1022
$ws.Range("E137").Value = 126

# Row 145 - Georgia
$ws.Range("B145").Value = 1006
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 883
$ws.Range("E145").Value = 108
